# Update the "想去人数" (interested-count) figures on the 展览 and 全部类型
# sheets to reflect the refreshed data pull.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5680
$ws1.Range("F7").Value = 157
$ws1.Range("F8").Value = 2561
$ws1.Range("F14").Value = 2398
$ws1.Range("F15").Value = 415

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5680
$ws4.Range("F9").Value = 157
$ws4.Range("F10").Value = 2561
$ws4.Range("F17").Value = 2398
$ws4.Range("F18").Value = 415
